$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 29
$ws.Range("I2").Value = 56
$ws.Range("J2").Value = 217
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 38
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 19
$ws.Range("T2").Value = 37
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 333
$ws.Range("X2").Value = 376
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 5
$ws.Range("AA2").Value = 1
